$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.977.45"
$ws.Range("E2").Value = "  +0.48%  "
$ws.Range("D3").Value = "1.653.85"
$ws.Range("E3").Value = "  +2.13%  "
$ws.Range("D4").Value = "'0.9992"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'308.84"
$ws.Range("E5").Value = "  +0.23%  "
$ws.Range("D6").Value = "'0.9994"
$ws.Range("E6").Value = "  -0.18%  "
$ws.Range("D7").Value = "'0.3906"
$ws.Range("E7").Value = "  -0.82%  "
$ws.Range("D8").Value = "'0.3833"
$ws.Range("E8").Value = "  -0.37%  "
$ws.Range("D9").Value = "'51.35"
$ws.Range("E9").Value = "  +3.25%  "
$ws.Range("D10").Value = "'1.355"
$ws.Range("E10").Value = "  -0.57%  "
$ws.Range("D11").Value = "'0.9999"
$ws.Range("E11").Value = "  -0.38%  "
$ws.Range("D12").Value = "'0.08446"
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("D13").Value = "'24.00"
$ws.Range("E13").Value = "  +0.69%  "
$ws.Range("D14").Value = "'7.125"
$ws.Range("E14").Value = "  +1.08%  "
$ws.Range("D15").Value = "'7.879"
$ws.Range("E15").Value = "  +4.09%  "
$ws.Range("D16").Value = "'0.00001314"
$ws.Range("E16").Value = "  +2.75%  "
$ws.Range("D17").Value = "1.651.51"
$ws.Range("E17").Value = "  +3.39%  "
$ws.Range("D18").Value = "'94.51"
$ws.Range("E18").Value = "  +0.59%  "
$ws.Range("E19").Value = "  +0.51%  "
$ws.Range("D20").Value = "'19.78"
$ws.Range("E20").Value = "  -1.50%  "
$ws.Range("D21").Value = "'6.906"
$ws.Range("E21").Value = "  +1.35%  "
$ws.Range("D22").Value = "'0.9990"
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("E23").Value = "  +1.69%  "
$ws.Range("D24").Value = "23.953.46"
$ws.Range("E24").Value = "  +0.38%  "
$ws.Range("D25").Value = "'2.483"
$ws.Range("E25").Value = "  +1.20%  "
$ws.Range("D26").Value = "'3.024"
$ws.Range("E26").Value = "  +5.82%  "
$ws.Range("D27").Value = "'22.04"
$ws.Range("E27").Value = "  -0.99%  "
$ws.Range("D28").Value = "'152.90"
$ws.Range("E28").Value = "  -2.33%  "
$ws.Range("D29").Value = "'5.444"
$ws.Range("E29").Value = "  +3.27%  "
$ws.Range("D30").Value = "'139.28"
$ws.Range("E30").Value = "  -0.62%  "
$ws.Range("D31").Value = "'7.746"
$ws.Range("E31").Value = "  -1.40%  "
$ws.Range("D32").Value = "'2.486"
$ws.Range("E32").Value = "  -0.97%  "
$ws.Range("D33").Value = "1.833.95"
$ws.Range("E33").Value = "  +2.25%  "
$ws.Range("D34").Value = "'1.038"
$ws.Range("E34").Value = "  +6.15%  "
$ws.Range("D35").Value = "'0.08132"
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("E36").Value = "  +2.91%  "
$ws.Range("D37").Value = "'6.755"
$ws.Range("E37").Value = "  +2.61%  "
$ws.Range("D38").Value = "'10.85"
$ws.Range("E38").Value = "  +4.98%  "
$ws.Range("D39").Value = "'0.2676"
$ws.Range("E39").Value = "  +0.36%  "
$ws.Range("D40").Value = "'0.09146"
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("D41").Value = "'0.7552"
$ws.Range("E41").Value = "  +0.54%  "
$ws.Range("D42").Value = "'13.51"
$ws.Range("E42").Value = "  -0.47%  "
$ws.Range("D43").Value = "'1.426"
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("D44").Value = "'16.27"
$ws.Range("E44").Value = "  +1.45%  "
$ws.Range("E45").Value = "  +0.71%  "
$ws.Range("D46").Value = "'2.457"
$ws.Range("E46").Value = "  -0.60%  "
$ws.Range("D47").Value = "'4.079"
$ws.Range("E47").Value = "  +0.19%  "
$ws.Range("D48").Value = "'0.9985"
$ws.Range("E48").Value = "  -0.19%  "
$ws.Range("D49").Value = "'0.08298"
$ws.Range("E49").Value = "  +0.76%  "
$ws.Range("D50").Value = "'134.50"
$ws.Range("E50").Value = "  +0.41%  "
$ws.Range("D51").Value = "'1.223"
$ws.Range("E51").Value = "  +0.82%  "

Write-Host "Updated cryptos list"
